# Mac/Safari compatibility & contrast fixes:
#  - mark the first "رافعات التحدي" entry (row 2) as Done
#  - record the new order placed 2025-05-19 (row 7) and mark it Done too

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# Row 2: flip "Done" flag to TRUE
$ws.Cells.Item(2, 9).Value = $true

# Row 7: new transaction
$row = 7

$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "رافعات التحدي"

$cellB = $ws.Cells.Item($row, 2)
$cellB.NumberFormat = "@"
$cellB.Value = "2025-05-19"

$cellC = $ws.Cells.Item($row, 3)
$cellC.NumberFormat = "@"
$cellC.Value = "مذكرات 100 حبة - شد 50 طباعة وجه واحد"

$cellD = $ws.Cells.Item($row, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "150"

$ws.Cells.Item($row, 5).Value = 500
$ws.Cells.Item($row, 6).Value = 75
$ws.Cells.Item($row, 7).Value = 575

# Column H (Actions) intentionally left blank for this row

$ws.Cells.Item($row, 9).Value = $true
